$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "41.540.17"
Set-TextValue "E2" "  +0.21%  "
Set-TextValue "D3" "2.460.47"
Set-TextValue "E3" "  +1.53%  "
Set-TextValue "E4" "  -1.57%  "
Set-TextValue "D5" "314.33"
Set-TextValue "E5" "  +2.24%  "
Set-TextValue "D6" "90.92"
Set-TextValue "E6" "  +3.78%  "
Set-TextValue "D7" "0.547"
Set-TextValue "E7" "  +4.26%  "
Set-TextValue "E8" "  -1.30%  "
Set-TextValue "D9" "0.510"
Set-TextValue "E9" "  +7.85%  "
Set-TextValue "D10" "32.41"
Set-TextValue "E10" "  +3.39%  "
Set-TextValue "D11" "0.0792"
Set-TextValue "E11" "  +5.27%  "
Set-TextValue "E12" "  +1.63%  "
Set-TextValue "D13" "2.841.21"
Set-TextValue "E13" "  +1.44%  "
Set-TextValue "D14" "6.83"
Set-TextValue "E14" "  +3.53%  "
Set-TextValue "D15" "15.77"
Set-TextValue "E15" "  +7.20%  "
Set-TextValue "D16" "2.434.44"
Set-TextValue "E16" "  -0.33%  "
Set-TextValue "D17" "0.770"
Set-TextValue "E17" "  +2.88%  "
Set-TextValue "D18" "41.531.88"
Set-TextValue "E18" "  +1.11%  "
Set-TextValue "D19" "6.47"
Set-TextValue "E19" "  +6.80%  "
Set-TextValue "D20" "0.0₃0935"
Set-TextValue "E20" "  +5.86%  "
Set-TextValue "D21" "70.84"
Set-TextValue "E21" "  +4.13%  "
Set-TextValue "D22" "11.26"
Set-TextValue "E22" "  +6.87%  "
Set-TextValue "D23" "236.76"
Set-TextValue "E23" "  +3.70%  "
Set-TextValue "D24" "2.72"
Set-TextValue "E24" "  +3.50%  "
Set-TextValue "E25" "  +0.00%  "
Set-TextValue "E26" "  +5.10%  "
Set-TextValue "D27" "24.20"
Set-TextValue "E27" "  +4.77%  "
Set-TextValue "E28" "  +2.76%  "
Set-TextValue "D29" "9.61"
Set-TextValue "E29" "  +3.11%  "
Set-TextValue "D30" "34.99"
Set-TextValue "E30" "  +3.23%  "
Set-TextValue "D31" "156.13"
Set-TextValue "E31" "  +3.91%  "
Set-TextValue "D32" "5.41"
Set-TextValue "E32" "  +4.62%  "
Set-TextValue "D33" "2.58"
Set-TextValue "E33" "  +1.78%  "
Set-TextValue "D34" "0.0755"
Set-TextValue "E34" "  +4.00%  "
Set-TextValue "D35" "17.24"
Set-TextValue "E35" "  +3.70%  "
Set-TextValue "D36" "2.40"
Set-TextValue "E36" "  -3.51%  "
Set-TextValue "E37" "  +0.41%  "
Set-TextValue "E38" "  +4.30%  "
Set-TextValue "D39" "0.102"
Set-TextValue "E39" "  +6.38%  "
Set-TextValue "D40" "1.76"
Set-TextValue "E40" "  +1.61%  "
Set-TextValue "D41" "3.95"
Set-TextValue "E41" "  +2.99%  "
Set-TextValue "E42" "  -1.93%  "
Set-TextValue "D43" "1.954.25"
Set-TextValue "E43" "  +2.55%  "
Set-TextValue "D44" "0.0281"
Set-TextValue "E44" "  +4.20%  "
Set-TextValue "D45" "18.63"
Set-TextValue "E45" "  -8.64%  "
Set-TextValue "D46" "2.89"
Set-TextValue "E46" "  +3.22%  "
Set-TextValue "E47" "  +6.64%  "
Set-TextValue "D48" "2.700.79"
Set-TextValue "E48" "  +0.54%  "
Set-TextValue "D49" "96.21"
Set-TextValue "E49" "  +4.41%  "
Set-TextValue "D50" "66.69"
Set-TextValue "E50" "  +4.84%  "
Set-TextValue "E51" "  +1.42%  "
